$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing data rows (rows 2-6) before writing the new table so
# no stale cells are left behind when the new table is shorter/longer.
$ws.Range("A2:E6").Clear()

$data = @(
    @('Hampton', 'Merrymen Cafe, 2 Small Street, Hampton VIC', '28/12/20 12:50pm-2:40pm', 'Case ate in store', 'new'),
    @('Hampton', 'Merrymen Cafe, 2 Small Street, Hampton VIC', '28/12/20 1:30pm-2:30pm', 'Case ate in store', 'old'),
    @('McKinnon', 'Hotlocks By Rachael Hairdresser, 260 McKinnon Road, McKinnon VIC 3204', '23/12/20 4:00pm-6:00pm', 'Case had hair cut in store', 'old'),
    @('Melbourne', 'Left Bank Melbourne Restaurant and Cocktail Bar, 1 Southbank Blvd', '25/12/20 12:00pm-3:00pm', 'Case attended bar', 'new'),
    @('Melbourne', 'Left Bank Melbourne, 1 Southbank Blvd', '25/12/20 12:00pm-3:00pm', 'Case ate in store', 'old'),
    @('Mordialloc', 'Woodlands Golf Club - club bar  109 White Street Mordialloc VIC 3195', '23/12/20 12:30pm-1:30pm', 'Case attended club house bar', 'old'),
    @('Southbank', 'Rockpool Bar and Grill, Crown Casino  8 Whiteman Street, Southbank', '23/12/20 8:00pm-11:00pm', 'Case ate in store', 'old')
)

$row = 2
foreach ($r in $data) {
    $ws.Range("A$row").Value = $r[0]
    $ws.Range("B$row").Value = $r[1]
    $ws.Range("C$row").Value = $r[2]
    $ws.Range("D$row").Value = $r[3]
    $ws.Range("E$row").Value = $r[4]
    $row++
}

# Re-fit column widths for the refreshed content (closest representable
# values given this runtime's column-width rounding granularity).
$ws.Columns.Item(1).ColumnWidth = 8.5
$ws.Columns.Item(2).ColumnWidth = 59
$ws.Columns.Item(3).ColumnWidth = 21.8333333333333
$ws.Columns.Item(4).ColumnWidth = 23.1666666666667

[void]$ws.Range("B13").Select()
